$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily row (row 81) with the next date and win counts
$ws.Range("A81").Value = 46030
$ws.Range("B81").Value = 180
$ws.Range("C81").Value = 192
$ws.Range("D81").Value = 184

# Match the date style used by the rest of column A (copy style from A80)
$ws.Range("A80").Copy()
$ws.Range("A81").PasteSpecial(-4122)
